# Update countries & provincias Spain
# Applies the COVID data refresh to the "Pais" worksheet:
#  - Polonia overtakes Rumania (rows 32/33 swap places, Polonia gets fresh data)
#  - Austria overtakes Armenia (rows 60/61 swap places, Austria gets fresh data)
#  - Hong Kong overtakes Botsuana (rows 130/131 swap places, Hong Kong gets fresh data)
#  - Several other countries get refreshed case numbers
#  - "Datos actualizados" timestamp bumped from 10:33 to 11:50

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 11:50"

# Helper: write a full data row (Pais, Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) using positional args
function Set-Row($Row, $Pais, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Row 4: Estados Unidos - refreshed numbers
Set-Row 4 "Estados Unidos" 8388013 214 5457912 2705369 0 2 224732

# Row 21: Alemania - refreshed numbers
Set-Row 21 "Alemania" 367419 438 291900 65652 0 1 9867

# Rows 32/33: Polonia moves ahead of Rumania
Set-Row 32 "Polonia" 183248 7482 94014 85620 0 41 3614
Set-Row 33 "Rumania" 180388 0 130894 43622 0 0 5872

# Row 39: Catar - refreshed numbers
Set-Row 39 "Catar" 129671 240 126650 2797 0 0 224

# Row 44: Oman - refreshed numbers
Set-Row 44 "Oman" 110594 641 96400 13080 0 13 1114

# Rows 60/61: Austria moves ahead of Armenia
Set-Row 60 "Austria" 65927 1121 50359 14664 0 11 904
Set-Row 61 "Armenia" 65460 766 48208 16161 0 10 1091

# Row 77: Afganistan - refreshed numbers
Set-Row 77 "Afganistan" 40287 87 33760 5030 0 5 1497

# Row 102: Finlandia - refreshed numbers
Set-Row 102 "Finlandia" 13555 131 9100 4104 0 0 351

# Row 107: Consejo Danes para los Refugiados - refreshed numbers
Set-Row 107 "Consejo Danes para los Refugiados" 11052 46 10357 392 0 1 303

# Rows 130/131: Hong Kong moves ahead of Botsuana
Set-Row 130 "Hong Kong" 5257 15 4982 170 0 0 105
Set-Row 131 "Botsuana" 5242 0 905 4317 0 0 20
